$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PROD")

# Update the value of A4 (was "jersey026.tt21.5") to the new token "dec9.examtaker.5"
$ws.Range("A4").Value = "dec9.examtaker.5"

# Move / update the active cell selection to A4
$ws.Activate()
$ws.Range("A4").Select()
